$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pages = @(
    "CategoryPages",
    "ComparePages",
    "DealerPages",
    "ErrorPages",
    "FCVPages",
    "HomeOffersPages",
    "LCertifiedPages",
    "MiscPages",
    "ModelPagesAccessories",
    "ModelPagesDesign",
    "ModelPagesFeatures",
    "ModelPagesGallery",
    "ModelPagesOffers",
    "ModelPagesOverview",
    "ModelPagesOwners",
    "ModelPagesPackages",
    "ModelPagesPerformance",
    "ModelPagesSafety",
    "ModelPagesSpecifications",
    "ModelPagesTechnology",
    "SponsoredAthletes"
)

$row = 4
foreach ($page in $pages) {
    $ws.Cells.Item($row, 1).Value = "test"
    $ws.Cells.Item($row, 2).Value = $page
    $row = $row + 1
}

$ws.Range("C26").Select()
